$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 96 (pushes existing rows 96..145 down to 97..146)
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new data record
$ws.Cells.Item(96, 1).Value = 5
$ws.Cells.Item(96, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(96, 3).Value = "Maule"
$ws.Cells.Item(96, 4).Value = 44460
$ws.Cells.Item(96, 5).Value = 7
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100108
$ws.Cells.Item(96, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(96, 9).Value = 100108005
$ws.Cells.Item(96, 10).Value = "Piña"
$ws.Cells.Item(96, 11).Value = "Caramelo"
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 54
$ws.Cells.Item(96, 14).Value = 19000
$ws.Cells.Item(96, 15).Value = 19000
$ws.Cells.Item(96, 16).Value = 19000
$ws.Cells.Item(96, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(96, 18).Value = "Ecuador"
$ws.Cells.Item(96, 19).Value = 1357
$ws.Cells.Item(96, 20).Value = 14
